# Updates cryptos list values (price & 1h volume change) and fixes the
# Filecoin/Stacks row ordering, per the 2024-07-11 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.969.39"
$ws.Range("E2").Value = "  -1.81%  "
# Row 3
$ws.Range("D3").Value = "3.099.46"
$ws.Range("E3").Value = "  -0.22%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").Value = "'525.50"
$ws.Range("E5").Value = "  +0.29%  "
# Row 6
$ws.Range("D6").Value = "'141.23"
$ws.Range("E6").Value = "  -1.92%  "
# Row 7
$ws.Range("E7").Value = "  +0.04%  "
# Row 8
$ws.Range("D8").Value = "3.098.52"
$ws.Range("E8").Value = "  -0.21%  "
# Row 9
$ws.Range("E9").Value = "  +0.99%  "
# Row 10
$ws.Range("E10").Value = "  -2.94%  "
# Row 11
$ws.Range("E11").Value = "  -1.22%  "
# Row 12
$ws.Range("E12").Value = "  +2.35%  "
# Row 13
$ws.Range("D13").Value = "3.630.21"
$ws.Range("E13").Value = "  -0.19%  "
# Row 14
$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = "  +2.75%  "
# Row 15
$ws.Range("D15").Value = "'25.62"
$ws.Range("E15").Value = "  -5.48%  "
# Row 16
$ws.Range("E16").Value = "  -1.35%  "
# Row 17
$ws.Range("D17").Value = "58.002.26"
$ws.Range("E17").Value = "  -1.67%  "
# Row 18
$ws.Range("D18").Value = "3.098.73"
$ws.Range("E18").Value = "  -0.18%  "
# Row 19
$ws.Range("D19").Value = "'6.09"
$ws.Range("E19").Value = "  -1.83%  "
# Row 20
$ws.Range("E20").Value = "  -2.19%  "
# Row 21
$ws.Range("D21").Value = "'7.97"
$ws.Range("E21").Value = "  -2.85%  "
# Row 22
$ws.Range("D22").Value = "'343.45"
$ws.Range("E22").Value = "  +0.09%  "
# Row 23
$ws.Range("E23").Value = "  -0.11%  "
# Row 24
$ws.Range("E24").Value = "  +0.68%  "
# Row 25
$ws.Range("D25").Value = "'67.39"
$ws.Range("E25").Value = "  +2.38%  "
# Row 26
$ws.Range("E26").Value = "  -0.75%  "
# Row 27
$ws.Range("E27").Value = "  +0.15%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("E28").Value = "  -1.30%  "
# Row 29
$ws.Range("E29").Value = "  +0.01%  "
# Row 30
$ws.Range("D30").Value = "'6.37"
$ws.Range("E30").Value = "  -5.27%  "
# Row 31
$ws.Range("D31").Value = "'7.28"
$ws.Range("E31").Value = "  -0.22%  "
# Row 32
$ws.Range("E32").Value = "  +1.15%  "
# Row 33
$ws.Range("D33").Value = "'20.96"
$ws.Range("E33").Value = "  -0.46%  "
# Row 34
$ws.Range("D34").Value = "'1.18"
$ws.Range("E34").Value = "  -2.59%  "
# Row 35
$ws.Range("D35").Value = "'158.81"
$ws.Range("E35").Value = "  +2.38%  "
# Row 36
$ws.Range("D36").Value = "'4.63"
$ws.Range("E36").Value = "  -0.35%  "
# Row 37
$ws.Range("D37").Value = "'6.16"
$ws.Range("E37").Value = "  -0.58%  "
# Row 38
$ws.Range("D38").Value = "'26.11"
$ws.Range("E38").Value = "  -3.01%  "
# Row 39
$ws.Range("E39").Value = "  -4.95%  "
# Row 40
$ws.Range("E40").Value = "  -2.70%  "
# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.58"
$ws.Range("E41").Value = "  +8.03%  "
# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.01"
$ws.Range("E42").Value = "  +1.51%  "
# Row 43
$ws.Range("D43").Value = "'0.684"
$ws.Range("E43").Value = "  +2.92%  "
# Row 44
$ws.Range("D44").Value = "3.137.93"
$ws.Range("E44").Value = "  -0.24%  "
# Row 45
$ws.Range("D45").Value = "'36.91"
$ws.Range("E45").Value = "  +0.09%  "
# Row 46
$ws.Range("E46").Value = "  +0.01%  "
# Row 47
$ws.Range("E47").Value = "  +1.81%  "
# Row 48
$ws.Range("D48").Value = "2.271.31"
$ws.Range("E48").Value = "  -0.82%  "
# Row 49
$ws.Range("D49").Value = "'0.993"
$ws.Range("E49").Value = "  +2.64%  "
# Row 50
$ws.Range("D50").Value = "'6.09"
$ws.Range("E50").Value = "  +1.13%  "
# Row 51
$ws.Range("D51").Value = "'20.44"
$ws.Range("E51").Value = "  -2.51%  "
